# Penalty Reward System (unfinished) - shift weekly forecast dates forward by
# one week and zero out / shrink the early MyForecast numbers, then refresh
# the dependent Summary statistics to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Forecast Comparison": shift Week_Start_Date (col B) forward by one
# week for every data row, and overwrite MyForecast (col D) with the new
# (much lower) forecast numbers.
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecastValues = @(0, 0, 0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $weekDates.Length; $i++) {
    $row = $i + 2

    # Keep the date text as plain text (it was stored as an inline string,
    # not a real date serial) by forcing the cell to Text format before
    # writing the value, so Excel doesn't silently reinterpret it.
    $dateCell = $wsForecast.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekDates[$i]

    $wsForecast.Cells.Item($row, 4).Value = $myForecastValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh the stats that depend on the forecast window.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

function Set-TextValue($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $wsSummary 2  2 "2022-12-25 to 2025-01-05"   # Historical Range
Set-TextValue $wsSummary 9  2 "17"                          # Total Forecast (16 Weeks)
Set-TextValue $wsSummary 10 2 "7"                           # Total Forecast (8 Weeks)
Set-TextValue $wsSummary 11 2 "2"                           # Total Forecast (4 Weeks)
Set-TextValue $wsSummary 12 2 "1"                           # Max Forecast
Set-TextValue $wsSummary 13 2 "2025-02-16"                  # Max Forecast Week
Set-TextValue $wsSummary 14 2 "0"                           # Min Forecast
Set-TextValue $wsSummary 15 2 "2025-01-12"                  # Min Forecast Week
